$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "nvTmX183"
$ws.Range("B2").Value = 23073105
$ws.Range("C2").Value = "jwwcmjx39"
$ws.Range("D2").Value = "kAX3m#9&"
$ws.Range("F2").Value = "JjmTLqRU"
$ws.Range("G2").Value = "xaTn"
